# Capacity Supply Curve.xlsx edit
# - Update the CSC-CSCSoCECBiaSY calibration block (B2:AE25) from 0.55 to 0.66
# - Move the active/selected tab from "About" to "CSC-CSCSoCECBiaSY"
# - Update the "About" sheet's lingering selection (was B12) to E19

$wb = $excel.ActiveWorkbook

# --- Update the capacity-factor calibration values (rows 2-25, cols B:AE) ---
$wsCalib = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$wsCalib.Range("B2:AE25").Value = 0.66

# --- Update the "About" sheet's remembered selection (no longer the active tab) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
$wsAbout.Range("E19").Select() | Out-Null

# --- Make CSC-CSCSoCECBiaSY the active / selected tab, matching the original selection ---
$wsCalib.Activate() | Out-Null
$wsCalib.Range("B2:AE25").Select() | Out-Null
